# Apply cryptos list update (auto-generated from canonical OOXML diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text would otherwise be mis-parsed as a number by Excel ---
# (single-dot decimal-looking strings). Force them to stay text: apply a temporary
# text number format, assign the literal string, then restore the default "Normal"
# style so the cell keeps no explicit style index (matching original formatting).
$textForcedCells = @('D5', 'D6', 'D11', 'D12', 'D14', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D29', 'D30', 'D31', 'D32', 'D33', 'D37', 'D40', 'D41', 'D42', 'D45', 'D47', 'D48', 'D50')
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '63.113.45'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '3.057.13'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '583.88'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').Value = '152.05'
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D9').Value = '3.058.27'
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('D11').Value = '5.85'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').Value = '0.449'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('E13').Value = '  -2.93%  '
$ws.Range('D14').Value = '36.21'
$ws.Range('E14').Value = '  -3.67%  '
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').Value = '3.557.12'
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '63.060.94'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').Value = '3.056.22'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('D20').Value = '481.81'
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range('D21').Value = '14.34'
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('D22').Value = '0.709'
$ws.Range('E22').Value = '  -1.56%  '
$ws.Range('D23').Value = '7.54'
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('D24').Value = '2.41'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').Value = '82.04'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('D26').Value = '12.70'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('E27').Value = '  +5.67%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = '7.43'
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '2.21'
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.66'
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').Value = '27.79'
$ws.Range('E33').Value = '  +1.81%  '
$ws.Range('E34').Value = '  -2.47%  '
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('D36').Value = '0.0₃0817'
$ws.Range('E36').Value = '  -4.05%  '
$ws.Range('D37').Value = '5.93'
$ws.Range('E37').Value = '  -3.29%  '
$ws.Range('E38').Value = '  -4.93%  '
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('D40').Value = '9.22'
$ws.Range('E40').Value = '  -1.40%  '
$ws.Range('D41').Value = '50.52'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('D42').Value = '429.13'
$ws.Range('E42').Value = '  -3.35%  '
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('E44').Value = '  +3.65%  '
$ws.Range('D45').Value = '0.0363'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('D46').Value = '2.846.02'
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('D47').Value = '38.17'
$ws.Range('E47').Value = '  -4.82%  '
$ws.Range('D48').Value = '127.28'
$ws.Range('E48').Value = '  -3.43%  '
$ws.Range('D50').Value = '25.18'
$ws.Range('E50').Value = '  -1.29%  '
$ws.Range('E51').Value = '  -1.15%  '

# Restore default styling on the text-forced cells (removes the temporary text format)
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).Style = "Normal"
}
